# Automatic update of files.
# - Column C ("Förändrad") on rows 2..26 bumped from 45252 to 45253 (2023-11-22 -> 2023-11-23)
# - Rows 27 and 28 (A 58731-2023 / A 58730-2023) removed
# - Row 26 loses its explicit custom row-height flag (reverts to default sizing)
# - Sheet dimension shrinks accordingly (A1:Y28 -> A1:Y26), which follows
#   automatically once the trailing rows are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date stamp for every data row that will remain.
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 45253
}

# Drop the last two rows (A 58731-2023, A 58730-2023). Deleting row 27 twice
# removes both, shifting nothing else since they were the last rows.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(27).Delete()

# Row 26 no longer carries an explicit customHeight flag after the trailing
# rows are removed; AutoFit clears the custom-height marker while leaving the
# effective height at the sheet default (15).
$ws.Rows.Item(26).AutoFit()
